$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 20,7

$arr[0,0] = "Arsenal"
$arr[0,1] = 4.76829268292683
$arr[0,2] = 8.755395683453237
$arr[0,3] = 0.6250581665891112
$arr[0,4] = 15
$arr[0,5] = 0.06095032755618211
$arr[0,6] = 49

$arr[1,0] = "Manchester City"
$arr[1,1] = 5.648741418764302
$arr[1,2] = 7.679319371727749
$arr[1,3] = 0.6848201774871555
$arr[1,4] = 22
$arr[1,5] = 0.04954533297673175
$arr[1,6] = 28

$arr[2,0] = "Aston Villa"
$arr[2,1] = 5.012121212121212
$arr[2,2] = 6.745019920318725
$arr[2,3] = 0.5270105379922352
$arr[2,4] = 29
$arr[2,5] = 0.06273731446323783
$arr[2,6] = -12

$arr[3,0] = "Manchester United"
$arr[3,1] = 4.339958158995816
$arr[3,2] = 7.149851632047477
$arr[3,3] = 0.5438832539846823
$arr[3,4] = 16
$arr[3,5] = 0.06906758756783424
$arr[3,6] = 23

$arr[4,0] = "Chelsea"
$arr[4,1] = 4.500515995872033
$arr[4,2] = 9.593519882179676
$arr[4,3] = 0.6126921818776578
$arr[4,4] = 22
$arr[4,5] = 0.05496314320475202
$arr[4,6] = 19

$arr[5,0] = "Liverpool"
$arr[5,1] = 4.878467635402906
$arr[5,2] = 7.691019786910198
$arr[5,3] = 0.6799582265261559
$arr[5,4] = 19
$arr[5,5] = 0.06433655133535339
$arr[5,6] = 16

$arr[6,0] = "Brentford"
$arr[6,1] = 5.085155350978136
$arr[6,2] = 7.107667210440456
$arr[6,3] = 0.4128239445187979
$arr[6,4] = 14
$arr[6,5] = 0.09452940572718875
$arr[6,6] = 22

$arr[7,0] = "Everton"
$arr[7,1] = 6.091482649842272
$arr[7,2] = 5.526785714285714
$arr[7,3] = 0.4536886577835554
$arr[7,4] = 26
$arr[7,5] = 0.07662062707934268
$arr[7,6] = -17

$arr[8,0] = "Bournemouth"
$arr[8,1] = 4.309803921568627
$arr[8,2] = 6.265402843601896
$arr[8,3] = 0.5041747683861375
$arr[8,4] = 14
$arr[8,5] = 0.09870303261491513
$arr[8,6] = 10

$arr[9,0] = "Newcastle United"
$arr[9,1] = 4.932182490752158
$arr[9,2] = 7.020618556701031
$arr[9,3] = 0.5413231874867893
$arr[9,4] = 26
$arr[9,5] = 0.0755036397494498
$arr[9,6] = 17

$arr[10,0] = "Sunderland"
$arr[10,1] = 5.429057888762769
$arr[10,2] = 8.075289575289576
$arr[10,3] = 0.3833404619332763
$arr[10,4] = 23
$arr[10,5] = 0.0869120654396728
$arr[10,6] = -13

$arr[11,0] = "Fulham"
$arr[11,1] = 5.413431269674711
$arr[11,2] = 9.52808988764045
$arr[11,3] = 0.4887622324981181
$arr[11,4] = 21
$arr[11,5] = 0.07578068555832063
$arr[11,6] = -27

$arr[12,0] = "Crystal Palace"
$arr[12,1] = 5.878419452887538
$arr[12,2] = 7.031353135313531
$arr[12,3] = 0.4300120048019208
$arr[12,4] = 19
$arr[12,5] = 0.09257244585105331
$arr[12,6] = 7

$arr[13,0] = "Brighton & Hove Albion"
$arr[13,1] = 4.715139442231076
$arr[13,2] = 8.073569482288828
$arr[13,3] = 0.5179259605278523
$arr[13,4] = 23
$arr[13,5] = 0.06755548329404114
$arr[13,6] = 1

$arr[14,0] = "Leeds United"
$arr[14,1] = 5.559793814432989
$arr[14,2] = 8.043261231281198
$arr[14,3] = 0.4071159767362299
$arr[14,4] = 13
$arr[14,5] = 0.07290772020230939
$arr[14,6] = 9

$arr[15,0] = "Tottenham Hotspur"
$arr[15,1] = 4.339698492462311
$arr[15,2] = 6.462482946793997
$arr[15,3] = 0.4862776752767528
$arr[15,4] = 31
$arr[15,5] = 0.07963386727688787
$arr[15,6] = -6

$arr[16,0] = "Nottingham Forest"
$arr[16,1] = 5.531521739130435
$arr[16,2] = 6.324812030075188
$arr[16,3] = 0.5122107969151671
$arr[16,4] = 27
$arr[16,5] = 0.08233023982357805
$arr[16,6] = -19

$arr[17,0] = "West Ham United"
$arr[17,1] = 5.457403651115619
$arr[17,2] = 6.289902280130293
$arr[17,3] = 0.3929440389294404
$arr[17,4] = 17
$arr[17,5] = 0.08742945100051308
$arr[17,6] = -31

$arr[18,0] = "Burnley"
$arr[18,1] = 5.986214209968186
$arr[18,2] = 6.547967479674797
$arr[18,3] = 0.3759959141981614
$arr[18,4] = 23
$arr[18,5] = 0.08957341169318067
$arr[18,6] = -44

$arr[19,0] = "Wolverhampton"
$arr[19,1] = 4.764212488350419
$arr[19,2] = 6.751196172248804
$arr[19,3] = 0.3748973727422003
$arr[19,4] = 24
$arr[19,5] = 0.08159739863834976
$arr[19,6] = -32

$ws.Range("A2:G21").Value = $arr
